$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.357.45"
$ws.Range("E2").Value = "  +1.92%  "
$ws.Range("D3").Value = "1.949.55"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9991"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4832"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2927"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06831"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "105.20"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("D12").Value = "1.962.80"
$ws.Range("E12").Value = "  +1.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07836"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.340"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6942"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "298.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +8.55%  "
$ws.Range("D17").Value = "31.349.86"
$ws.Range("E17").Value = "  +1.90%  "
$ws.Range("D18").Value = "2.226.23"
$ws.Range("E18").Value = "  +1.84%  "
$ws.Range("D19").Value = "13.08"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007644"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.606"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9989"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9991"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.493"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.628"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.30%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "169.06"
$ws.Range("E26").Value = "  +2.51%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "19.93"
$ws.Range("E27").Value = "  +2.48%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.151"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "1.401"
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("E30").Value = "  -1.75%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "4.657"
$ws.Range("E31").Value = "  +2.05%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "1.538"
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "4.383"
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04856"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7484"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.140"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.737"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01971"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.671"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.33%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.653"
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "77.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "2.044"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8764"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4399"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "106.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9992"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.028.36"
$ws.Range("E47").Value = "  +4.52%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.636"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.22%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.189"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1220"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.22%  "
